$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 221-223: add weekly-average columns (H, J, M, O, Q) ---
$ws.Range("H221").Value = 10.28571428571429
$ws.Range("J221").Value = 5.288307094181478
$ws.Range("M221").Value = 30.71428571428572
$ws.Range("O221").Value = 0.4285714285714285
$ws.Range("Q221").Value = 20

$ws.Range("H222").Value = 10.28571428571429
$ws.Range("J222").Value = 5.288307094181478
$ws.Range("M222").Value = 30.71428571428572
$ws.Range("O222").Value = 0.4285714285714285
$ws.Range("Q222").Value = 20

$ws.Range("H223").Value = 18.42857142857143
$ws.Range("J223").Value = 7.960157178854959
$ws.Range("M223").Value = 40.57142857142857
$ws.Range("O223").Value = 0.4285714285714285
$ws.Range("Q223").Value = 21.71428571428572

# --- Append new rows 224-238 ---
# Copy the date-cell number format from A223 so new date cells (A224:A238) match style s="2"
$ws.Range("A223").Copy()
$ws.Range("A224:A238").PasteSpecial(-4122)

# Row 224
$ws.Range("A224").Value = 44113
$ws.Range("B224").Value = 3027
$ws.Range("C224").Value = 1515
$ws.Range("D224").Value = 109
$ws.Range("E224").Value = 2664
$ws.Range("F224").Value = 254
$ws.Range("G224").Value = 11
$ws.Range("H224").Value = 18.42857142857143
$ws.Range("I224").Value = 4.526748971193416
$ws.Range("J224").Value = 7.960157178854959
$ws.Range("K224").Value = 39
$ws.Range("L224").Value = 19
$ws.Range("M224").Value = 40.57142857142857
$ws.Range("N224").Value = 2
$ws.Range("O224").Value = 0.4285714285714285
$ws.Range("P224").Value = 26
$ws.Range("Q224").Value = 21.71428571428572
$ws.Range("R224").Value = 543.4216715198213
$ws.Range("S224").Value = 608.5314910025708
$ws.Range("T224").Value = 48.20051413881748
$ws.Range("U224").Value = 38.59783923910194

# Row 225
$ws.Range("A225").Value = 44114
$ws.Range("B225").Value = 3027
$ws.Range("C225").Value = 1515
$ws.Range("D225").Value = 109
$ws.Range("E225").Value = 2664
$ws.Range("F225").Value = 254
$ws.Range("G225").Value = 0
$ws.Range("H225").Value = 17.71428571428572
$ws.Range("I225").Value = 0
$ws.Range("J225").Value = 8.361134262921452
$ws.Range("K225").Value = 0
$ws.Range("L225").Value = 0
$ws.Range("M225").Value = 42.42857142857143
$ws.Range("N225").Value = 0
$ws.Range("O225").Value = 0.2857142857142857
$ws.Range("P225").Value = 0
$ws.Range("Q225").Value = 24.42857142857143
$ws.Range("R225").Value = 543.4216715198213
$ws.Range("S225").Value = 608.5314910025708
$ws.Range("T225").Value = 48.20051413881748
$ws.Range("U225").Value = 38.59783923910194

# Row 226
$ws.Range("A226").Value = 44115
$ws.Range("B226").Value = 3167
$ws.Range("C226").Value = 1585
$ws.Range("D226").Value = 109
$ws.Range("E226").Value = 2723
$ws.Range("F226").Value = 335
$ws.Range("G226").Value = 81
$ws.Range("H226").Value = 18.57142857142857
$ws.Range("I226").Value = 31.88976377952756
$ws.Range("J226").Value = 7.72594239581844
$ws.Range("K226").Value = 140
$ws.Range("L226").Value = 70
$ws.Range("M226").Value = 45.85714285714285
$ws.Range("N226").Value = 0
$ws.Range("O226").Value = 0.2857142857142857
$ws.Range("P226").Value = 59
$ws.Range("Q226").Value = 27
$ws.Range("R226").Value = 568.5551482336551
$ws.Range("S226").Value = 636.6484575835475
$ws.Range("T226").Value = 58.6439588688946
$ws.Range("U226").Value = 50.98505276234862

# Row 227
$ws.Range("A227").Value = 44116
$ws.Range("B227").Value = 3167
$ws.Range("C227").Value = 1585
$ws.Range("D227").Value = 109
$ws.Range("E227").Value = 2723
$ws.Range("F227").Value = 335
$ws.Range("G227").Value = 0
$ws.Range("H227").Value = 28
$ws.Range("I227").Value = 0
$ws.Range("J227").Value = 9.688789765815759
$ws.Range("K227").Value = 0
$ws.Range("L227").Value = 0
$ws.Range("M227").Value = 57.85714285714285
$ws.Range("N227").Value = 0
$ws.Range("O227").Value = 0.2857142857142857
$ws.Range("P227").Value = 0
$ws.Range("Q227").Value = 29.57142857142857
$ws.Range("R227").Value = 568.5551482336551
$ws.Range("S227").Value = 636.6484575835475
$ws.Range("T227").Value = 58.6439588688946
$ws.Range("U227").Value = 50.98505276234862

# Row 228
$ws.Range("A228").Value = 44117
$ws.Range("B228").Value = 3188
$ws.Range("C228").Value = 1594
$ws.Range("D228").Value = 109
$ws.Range("E228").Value = 2772
$ws.Range("F228").Value = 307
$ws.Range("G228").Value = -28
$ws.Range("H228").Value = 36.42857142857143
$ws.Range("I228").Value = -8.358208955223882
$ws.Range("J228").Value = 11.32001566937778
$ws.Range("K228").Value = 21
$ws.Range("L228").Value = 9
$ws.Range("M228").Value = 68.85714285714286
$ws.Range("N228").Value = 0
$ws.Range("O228").Value = 0
$ws.Range("P228").Value = 49
$ws.Range("Q228").Value = 32.42857142857143
$ws.Range("R228").Value = 572.3251697407302
$ws.Range("S228").Value = 640.2634961439588
$ws.Range("T228").Value = 59.84897172236504
$ws.Range("U228").Value = 53.31887560006176

# Row 229
$ws.Range("A229").Value = 44118
$ws.Range("B229").Value = 3257
$ws.Range("C229").Value = 1626
$ws.Range("D229").Value = 109
$ws.Range("E229").Value = 2806
$ws.Range("F229").Value = 342
$ws.Range("G229").Value = 35
$ws.Range("H229").Value = 36.42857142857143
$ws.Range("I229").Value = 11.40065146579805
$ws.Range("J229").Value = 11.32001566937778
$ws.Range("K229").Value = 69
$ws.Range("L229").Value = 32
$ws.Range("M229").Value = 68.85714285714286
$ws.Range("N229").Value = 0
$ws.Range("O229").Value = 0
$ws.Range("P229").Value = 34
$ws.Range("Q229").Value = 32.42857142857143
$ws.Range("R229").Value = 584.7123832639769
$ws.Range("S229").Value = 653.1169665809769
$ws.Range("T229").Value = 63.46401028277634
$ws.Range("U229").Value = 57.62747160814756

# Row 230
$ws.Range("A230").Value = 44119
$ws.Range("B230").Value = 3393
$ws.Range("C230").Value = 1672
$ws.Range("D230").Value = 109
$ws.Range("E230").Value = 2845
$ws.Range("F230").Value = 439
$ws.Range("G230").Value = 97
$ws.Range("H230").Value = 45.28571428571428
$ws.Range("I230").Value = 28.3625730994152
$ws.Range("J230").Value = 10.77780692287777
$ws.Range("K230").Value = 136
$ws.Range("L230").Value = 46
$ws.Range("M230").Value = 87.28571428571429
$ws.Range("N230").Value = 0
$ws.Range("O230").Value = 0.1428571428571428
$ws.Range("P230").Value = 39
$ws.Range("Q230").Value = 41.85714285714285
$ws.Range("R230").Value = 609.1277606431297
$ws.Range("S230").Value = 671.5938303341902
$ws.Range("T230").Value = 70.69408740359897
$ws.Range("U230").Value = 72.70755763644785

# Row 231
$ws.Range("A231").Value = 44120
$ws.Range("B231").Value = 3509
$ws.Range("C231").Value = 1729
$ws.Range("D231").Value = 109
$ws.Range("E231").Value = 2891
$ws.Range("F231").Value = 509
$ws.Range("G231").Value = 70
$ws.Range("H231").Value = 45.28571428571428
$ws.Range("I231").Value = 15.94533029612756
$ws.Range("J231").Value = 10.77780692287777
$ws.Range("K231").Value = 116
$ws.Range("L231").Value = 57
$ws.Range("M231").Value = 87.28571428571429
$ws.Range("N231").Value = 0
$ws.Range("O231").Value = 0.1428571428571428
$ws.Range("P231").Value = 46
$ws.Range("Q231").Value = 41.85714285714285
$ws.Range("R231").Value = 629.9526413488778
$ws.Range("S231").Value = 694.4890745501285
$ws.Range("T231").Value = 85.95758354755783
$ws.Range("U231").Value = 86.53096982905646

# Row 232
$ws.Range("A232").Value = 44121
$ws.Range("B232").Value = 3509
$ws.Range("C232").Value = 1729
$ws.Range("D232").Value = 109
$ws.Range("E232").Value = 2891
$ws.Range("F232").Value = 509
$ws.Range("G232").Value = 0
$ws.Range("H232").Value = 45
$ws.Range("I232").Value = 0
$ws.Range("J232").Value = 11.31451863164333
$ws.Range("K232").Value = 0
$ws.Range("L232").Value = 0
$ws.Range("M232").Value = 89.71428571428571
$ws.Range("N232").Value = 0
$ws.Range("O232").Value = 0.1428571428571428
$ws.Range("P232").Value = 0
$ws.Range("Q232").Value = 44.57142857142857
$ws.Range("R232").Value = 629.9526413488778
$ws.Range("S232").Value = 694.4890745501285
$ws.Range("T232").Value = 85.95758354755783
$ws.Range("U232").Value = 86.53096982905646

# Row 233
$ws.Range("A233").Value = 44122
$ws.Range("B233").Value = 3778
$ws.Range("C233").Value = 1838
$ws.Range("D233").Value = 110
$ws.Range("E233").Value = 3016
$ws.Range("F233").Value = 652
$ws.Range("G233").Value = 143
$ws.Range("H233").Value = 43.28571428571428
$ws.Range("I233").Value = 28.09430255402751
$ws.Range("J233").Value = 10.21410402169239
$ws.Range("K233").Value = 269
$ws.Range("L233").Value = 109
$ws.Range("M233").Value = 93.57142857142857
$ws.Range("N233").Value = 1
$ws.Range("O233").Value = 0.1428571428571428
$ws.Range("P233").Value = 125
$ws.Range("Q233").Value = 50.14285714285715
$ws.Range("R233").Value = 678.2448216061728
$ws.Range("S233").Value = 738.271208226221
$ws.Range("T233").Value = 101.6227506426735
$ws.Range("U233").Value = 109.6896733725176

# Row 234
$ws.Range("A234").Value = 44123
$ws.Range("B234").Value = 3778
$ws.Range("C234").Value = 1838
$ws.Range("D234").Value = 110
$ws.Range("E234").Value = 3016
$ws.Range("F234").Value = 652
$ws.Range("G234").Value = 0
$ws.Range("H234").Value = 45.14285714285715
$ws.Range("I234").Value = 0
$ws.Range("J234").Value = 8.598631231188998
$ws.Range("K234").Value = 0
$ws.Range("L234").Value = 0
$ws.Range("M234").Value = 101.4285714285714
$ws.Range("N234").Value = 0
$ws.Range("O234").Value = 0.2857142857142857
$ws.Range("P234").Value = 0
$ws.Range("Q234").Value = 56
$ws.Range("R234").Value = 678.2448216061728
$ws.Range("S234").Value = 738.271208226221
$ws.Range("T234").Value = 101.6227506426735
$ws.Range("U234").Value = 109.6896733725176

# Row 235
$ws.Range("A235").Value = 44124
$ws.Range("B235").Value = 3816
$ws.Range("C235").Value = 1855
$ws.Range("D235").Value = 110
$ws.Range("E235").Value = 3084
$ws.Range("F235").Value = 622
$ws.Range("G235").Value = -30
$ws.Range("H235").Value = 48.28571428571428
$ws.Range("I235").Value = -4.601226993865031
$ws.Range("J235").Value = 8.061502683681653
$ws.Range("K235").Value = 38
$ws.Range("L235").Value = 17
$ws.Range("M235").Value = 111.2857142857143
$ws.Range("N235").Value = 0
$ws.Range("O235").Value = 0.2857142857142857
$ws.Range("P235").Value = 68
$ws.Range("Q235").Value = 62.71428571428572
$ws.Range("R235").Value = 685.066765285642
$ws.Range("S235").Value = 745.0996143958869
$ws.Range("T235").Value = 104.836118251928
$ws.Range("U235").Value = 112.7415955449117

# Row 236
$ws.Range("A236").Value = 44125
$ws.Range("B236").Value = 3912
$ws.Range("C236").Value = 1894
$ws.Range("D236").Value = 110
$ws.Range("E236").Value = 3157
$ws.Range("F236").Value = 645
$ws.Range("G236").Value = 23
$ws.Range("I236").Value = 3.697749196141479
$ws.Range("K236").Value = 96
$ws.Range("L236").Value = 39
$ws.Range("N236").Value = 0
$ws.Range("P236").Value = 73
$ws.Range("R236").Value = 702.3011493179852
$ws.Range("S236").Value = 760.7647814910026
$ws.Range("T236").Value = 107.6478149100257
$ws.Range("U236").Value = 117.5887660540082

# Row 237
$ws.Range("A237").Value = 44126
$ws.Range("B237").Value = 4103
$ws.Range("C237").Value = 1973
$ws.Range("D237").Value = 111
$ws.Range("E237").Value = 3237
$ws.Range("F237").Value = 755
$ws.Range("G237").Value = 110
$ws.Range("I237").Value = 17.05426356589147
$ws.Range("K237").Value = 191
$ws.Range("L237").Value = 79
$ws.Range("N237").Value = 1
$ws.Range("P237").Value = 80
$ws.Range("R237").Value = 736.5903925490013
$ws.Range("S237").Value = 792.4967866323908
$ws.Range("T237").Value = 120.9029562982005
$ws.Range("U237").Value = 127.4626319058715

# Row 238
$ws.Range("A238").Value = 44127
$ws.Range("B238").Value = 4288
$ws.Range("C238").Value = 2031
$ws.Range("D238").Value = 111
$ws.Range("E238").Value = 3330
$ws.Range("F238").Value = 847
$ws.Range("G238").Value = 92
$ws.Range("I238").Value = 12.18543046357616
$ws.Range("K238").Value = 185
$ws.Range("L238").Value = 58
$ws.Range("N238").Value = 0
$ws.Range("P238").Value = 93
$ws.Range("R238").Value = 769.802486777996
$ws.Range("S238").Value = 815.7937017994858
$ws.Range("T238").Value = 121.3046272493573
$ws.Range("U238").Value = 139.8498454291182
